$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.611300000000004
$ws.Range("B21").Value = 9.4297
$ws.Range("B23").Value = 8.869799999999996
$ws.Range("C24").Value = -13.57999999999999
$ws.Range("B25").Value = 6.081899999999996
$ws.Range("C28").Value = -13.7972
$ws.Range("C36").Value = -11.9264
$ws.Range("C45").Value = -13.67399999999999
$ws.Range("C48").Value = -12.16709999999999
$ws.Range("C49").Value = -13.7908
$ws.Range("C52").Value = -10.7776
$ws.Range("B53").Value = 6.321299999999993
$ws.Range("C53").Value = -10.87420000000001
$ws.Range("C54").Value = -13.38480000000001
$ws.Range("B57").Value = 4.950199999999993
$ws.Range("B59").Value = 5.200899999999999
$ws.Range("B69").Value = 5.251499999999994
$ws.Range("C70").Value = -12.2576
$ws.Range("B79").Value = 9.859200000000008
$ws.Range("B83").Value = 5.668199999999997
$ws.Range("C86").Value = -13.8331
$ws.Range("C87").Value = -12.98679999999999
$ws.Range("B93").Value = 5.6311
$ws.Range("C101").Value = -13.0037
